$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)
$dataSheet = $wb.Worksheets.Item("data")

# --- Add the new "feasibility" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "feasibility"

# The new sheet mirrors the "Manufacturability" column (E) of the "data"
# sheet, but stores it as binary feasibility flags (1/0) instead of the
# v/XX text, reusing the same header/format/conditional fill styling.

# Row 1: header "Manufacturability"
$dataSheet.Range("E1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Value = "Manufacturability"

# Row 2: units placeholder "--" (quote-prefixed, like the source cell)
$dataSheet.Range("E2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Value = "'--"

# Rows 3-22: feasibility flag per run, same fills as the source column
$values = @(1, 0, 0, 0, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 0, 1, 1, 0, 1, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 3
    $srcCell = "E$row"
    $dstCell = "A$row"
    $dataSheet.Range($srcCell).Copy()
    $ws.Range($dstCell).PasteSpecial(-4122)
    $ws.Range($dstCell).Value = $values[$i]
}

# Match source column width (Manufacturability column on "data")
$ws.Columns("A:A").ColumnWidth = $dataSheet.Columns("E:E").ColumnWidth

# --- Selections / active sheet bookkeeping (mirrors the recorded session) ---
$sheet1.Select() | Out-Null
$sheet1.Range("A1:I22").Select() | Out-Null

$dataSheet.Select() | Out-Null
$dataSheet.Range("E1:E1048576").Select() | Out-Null

$ws.Select() | Out-Null
$ws.Range("A23").Select() | Out-Null
